# Loan RBI, Variable Instalments
#
# The "Repayment schedule" sheet gets a new (blank) column inserted
# before the existing "Late" column so that a new "Outstanding"-style
# column can sit between "In Advance" and "Late"/"Outstanding"/"Original"
# for variable instalment tracking. All the following columns shift one
# place to the right (N->O, O->P, P->Q) and the sheet becomes the
# active / selected tab of the workbook (instead of "Transactions").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at N (pushes old N,O,P to O,P,Q)
$ws.Columns("N:N").Insert()

# Give the freshly inserted column the same on-screen width as the
# other data columns in this table (~11 characters)
$ws.Columns("N:N").ColumnWidth = 10.14

# Make "Repayment schedule" the active sheet / tab, and move the
# selection to L19, matching the edited workbook's saved view state
$ws.Activate()
$ws.Range("L19").Select()
